$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1: rows 102-104 (flash-file-overlay) - set string columns in C,D,I column-major order first
$ws.Range("C102").Value = "flash-file-overlay"
$ws.Range("C103").Value = "flash-file-overlay"
$ws.Range("C104").Value = "flash-file-overlay"
$ws.Range("D102").Value = "[suspicious-conversation-email, ad-click,suspicious-link,credentials-theft,frequents-spams]"
$ws.Range("D103").Value = "[suspicious-conversation-email, ad-click,suspicious-link,credentials-theft,frequents-spams]"
$ws.Range("D104").Value = "[suspicious-conversation-email, ad-click,suspicious-link,credentials-theft,frequents-spams]"
$ws.Range("I102").Value = "support_invisible_flash"
$ws.Range("I103").Value = "support_invisible_flash"
$ws.Range("I104").Value = "support_invisible_flash"
# Block 1: remaining columns
$ws.Range("A102").Value = "facebook"
$ws.Range("F102").Value = "technology"
$ws.Range("G102").Value = "north-america"
$ws.Range("H102").Value = "high"
$ws.Range("J102").Value = "[confidentiality:read_data;integrity:modify_data]"
$ws.Range("K102").Value = "medium"
$ws.Range("L102").Value = "medium"
$ws.Range("A103").Value = "alibaba"
$ws.Range("F103").Value = "technology"
$ws.Range("G103").Value = "asia"
$ws.Range("H103").Value = "high"
$ws.Range("J103").Value = "[confidentiality:read_data;integrity:modify_data]"
$ws.Range("K103").Value = "medium"
$ws.Range("L103").Value = "medium"
$ws.Range("A104").Value = "credit_agricole"
$ws.Range("F104").Value = "finance"
$ws.Range("G104").Value = "europe"
$ws.Range("H104").Value = "high"
$ws.Range("J104").Value = "[confidentiality:read_data;integrity:modify_data]"
$ws.Range("K104").Value = "medium"
$ws.Range("L104").Value = "medium"
$ws.Range("B102").Value = 43444
$ws.Range("E102").Value = 2000
$ws.Range("B103").Value = 44251
$ws.Range("E103").Value = 15000
$ws.Range("B104").Value = 44092
$ws.Range("E104").Value = 220

# Block 2: rows 105-110 (android-activity-hijack + tapjacking) - column-major C,D,I across both sub-groups
$ws.Range("C105").Value = "android-activity-hijack"
$ws.Range("C106").Value = "android-activity-hijack"
$ws.Range("C107").Value = "android-activity-hijack"
$ws.Range("C108").Value = "tapjacking"
$ws.Range("C109").Value = "tapjacking"
$ws.Range("C110").Value = "tapjacking"
$ws.Range("D105").Value = "[app-download,credentials-theft,suspicious-code-modifications]"
$ws.Range("D106").Value = "[app-download,credentials-theft,suspicious-code-modifications]"
$ws.Range("D107").Value = "[app-download,credentials-theft,suspicious-code-modifications]"
$ws.Range("D108").Value = "[pop-up-windows,app-download,ad-click,frequents-spams]"
$ws.Range("D109").Value = "[pop-up-windows,app-download,ad-click,frequents-spams]"
$ws.Range("D110").Value = "[pop-up-windows,app-download,ad-click,frequents-spams]"
$ws.Range("I105").Value = "previously_installed_malicious_application_on_android_device"
$ws.Range("I106").Value = "previously_installed_malicious_application_on_android_device"
$ws.Range("I107").Value = "previously_installed_malicious_application_on_android_device"
$ws.Range("I108").Value = "previously_installed_malicious_application"
$ws.Range("I109").Value = "previously_installed_malicious_application"
$ws.Range("I110").Value = "previously_installed_malicious_application"
# Block 2: remaining columns
$ws.Range("A105").Value = "samsung"
$ws.Range("F105").Value = "technology"
$ws.Range("G105").Value = "asia"
$ws.Range("H105").Value = "high"
$ws.Range("J105").Value = "[confidentiality:read_data]"
$ws.Range("K105").Value = "medium"
$ws.Range("L105").Value = "medium"
$ws.Range("A106").Value = "Intesa"
$ws.Range("F106").Value = "finance"
$ws.Range("G106").Value = "europe"
$ws.Range("H106").Value = "high"
$ws.Range("J106").Value = "[confidentiality:read_data]"
$ws.Range("K106").Value = "medium"
$ws.Range("L106").Value = "medium"
$ws.Range("A107").Value = "Erste bank"
$ws.Range("F107").Value = "finance"
$ws.Range("G107").Value = "europe"
$ws.Range("H107").Value = "high"
$ws.Range("J107").Value = "[confidentiality:read_data]"
$ws.Range("K107").Value = "medium"
$ws.Range("L107").Value = "medium"
$ws.Range("A108").Value = "credit_agricole"
$ws.Range("F108").Value = "finance"
$ws.Range("G108").Value = "europe"
$ws.Range("H108").Value = "high"
$ws.Range("J108").Value = "[confidentiality:other]"
$ws.Range("K108").Value = "low"
$ws.Range("L108").Value = "low"
$ws.Range("A109").Value = "credit_agricole"
$ws.Range("F109").Value = "finance"
$ws.Range("G109").Value = "europe"
$ws.Range("H109").Value = "high"
$ws.Range("J109").Value = "[confidentiality:other]"
$ws.Range("K109").Value = "low"
$ws.Range("L109").Value = "low"
$ws.Range("A110").Value = "samsung"
$ws.Range("F110").Value = "technology"
$ws.Range("G110").Value = "asia"
$ws.Range("H110").Value = "high"
$ws.Range("J110").Value = "[confidentiality:other]"
$ws.Range("K110").Value = "low"
$ws.Range("L110").Value = "low"
$ws.Range("B105").Value = 43824
$ws.Range("E105").Value = 1500
$ws.Range("B106").Value = 44321
$ws.Range("E106").Value = 250
$ws.Range("B107").Value = 43936
$ws.Range("E107").Value = 175
$ws.Range("B108").Value = 44147
$ws.Range("E108").Value = 100
$ws.Range("B109").Value = 44056
$ws.Range("E109").Value = 100
$ws.Range("B110").Value = 43308
$ws.Range("E110").Value = 1500

# Block 3: rows 111-113 (rooting-sim-cards) - brand new rows; B/E need style copied from an existing
# styled row (col B style s=4 date fmt, col E style s=6 thousands fmt) since these are new cells
# with no pre-existing style to inherit.
$ws.Range("B90").Copy()
$ws.Range("B111:B113").PasteSpecial(-4122)
$ws.Range("E90").Copy()
$ws.Range("E111:E113").PasteSpecial(-4122)

# A,C,D,I column-major first (A111/A113 introduce "iPhone")
$ws.Range("A111").Value = "iPhone"
$ws.Range("A112").Value = "samsung"
$ws.Range("A113").Value = "iPhone"
$ws.Range("C111").Value = "rooting-sim-cards"
$ws.Range("C112").Value = "rooting-sim-cards"
$ws.Range("C113").Value = "rooting-sim-cards"
$ws.Range("D111").Value = "[update,suspicious-conversation-phone-message,credentials-theft]"
$ws.Range("D112").Value = "[update,suspicious-conversation-phone-message,credentials-theft]"
$ws.Range("D113").Value = "[update,suspicious-conversation-phone-message,credentials-theft]"
$ws.Range("I111").Value = "sim-card-that-relies-on-des-cipher"
$ws.Range("I112").Value = "sim-card-that-relies-on-des-cipher"
$ws.Range("I113").Value = "sim-card-that-relies-on-des-cipher"
# Block 3: remaining columns
$ws.Range("F111").Value = "technology"
$ws.Range("G111").Value = "north-america"
$ws.Range("H111").Value = "medium"
$ws.Range("J111").Value = "[confidentiality:execute_unauthorized_commands]"
$ws.Range("K111").Value = "high"
$ws.Range("L111").Value = "high"
$ws.Range("F112").Value = "technology"
$ws.Range("G112").Value = "asia"
$ws.Range("H112").Value = "medium"
$ws.Range("J112").Value = "[confidentiality:execute_unauthorized_commands]"
$ws.Range("K112").Value = "high"
$ws.Range("L112").Value = "high"
$ws.Range("F113").Value = "technology"
$ws.Range("G113").Value = "north-america"
$ws.Range("H113").Value = "medium"
$ws.Range("J113").Value = "[confidentiality:execute_unauthorized_commands]"
$ws.Range("K113").Value = "high"
$ws.Range("L113").Value = "high"
$ws.Range("B111").Value = 43678
$ws.Range("E111").Value = 4200
$ws.Range("B112").Value = 43168
$ws.Range("E112").Value = 1500
$ws.Range("B113").Value = 44223
$ws.Range("E113").Value = 4200

$ws.Range("K113").Select()
Write-Host "done"